# Updates the "cryptos" price list: refreshed Price/Volume(1h) figures for
# most rows, plus a few coins that moved rank (Stacks/OKB swap rows 33/34,
# PEPE/Monero swap rows 46/47) and one coin replaced outright (ThetaToken ->
# CoreDAO in row 51).
#
# All of B/C/D/E are plain text columns in the source sheet (no numeric
# cells), so every assignment below is written with a leading apostrophe to
# force Excel to store it as text rather than auto-coercing numeric-looking
# strings (e.g. "535.41", "56.86") into real floating point numbers. The
# Style reset afterwards clears the "quote prefix" formatting flag that the
# apostrophe entry leaves behind, so the cell itself keeps its original
# (default) style - only the value's type/content changes, matching the
# source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.886.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.74%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.062.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.25%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'535.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.71%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'132.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.31%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.053.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.29%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.491"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.24%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -4.68%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'6.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -8.84%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.68%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.33%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'33.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.07%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.556.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.19%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'62.875.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.77%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.61%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.064.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.10%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'6.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.59%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'479.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.02%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.64%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.38%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -4.19%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'78.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.64%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'11.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.77%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.31%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.32%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.56%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.10%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'25.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.16%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -9.09%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.85%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'OKB"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'56.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.14%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'Stacks"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'2.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.53%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.23%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.16%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'474.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -12.00%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'3.084.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.46%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.0390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.29%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.0788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.70%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -3.21%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'8.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.51%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.96%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.250"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.74%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Monero"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'121.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.16%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'PEPE"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0₃0531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.71%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.70%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'24.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.11%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.18%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'CoreDAO"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'2.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.57%  "
$ws.Range("E51").Style = "Normal"
